# The sheet was wiped back to an empty template: all header/data cell
# values (A1:B6 - "Item"/"Price" headers plus the Widget rows) are removed
# while the existing cell styles/formatting stay untouched, and the
# selected cell moves to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1:B6").ClearContents()

$ws.Range("C6").Select()
